$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text change ---
$ws.Range("B1").Value = "Value (g)"

# --- Column A width (bestFit) ---
$ws.Columns("A").ColumnWidth = 17.6640625

# --- Selection change ---
$ws.Range("D8").Select()

# --- Helper cell used to coerce plain numeric-looking strings into real
#     text cells (apostrophe-prefix forces Excel to store them verbatim as
#     text instead of re-parsing them back into numbers); PasteSpecial
#     Values copies only the resulting text into the destination cell,
#     leaving the destination's own formatting (default/general) alone.
$helper = $ws.Range("Z1")

function Set-TextValue($addr, $text) {
    $helper.Value = "'" + $text
    $helper.Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

Set-TextValue "B2" "120.92"
Set-TextValue "B3" "12.65"
Set-TextValue "B4" "46.18"
Set-TextValue "B5" "21.01"

# Row 6 - parsed as a thousands-grouped number (decimal marker lost), scaled x10000
$ws.Range("B6").Value = "409,275"

Set-TextValue "B7" "120.92"
Set-TextValue "B8" "44.08"
Set-TextValue "B9" "90.3"
Set-TextValue "B10" "90.3"
Set-TextValue "B11" "60.46"
Set-TextValue "B12" "12.65"
Set-TextValue "B13" "133.44"
Set-TextValue "B14" "10.93"

# Row 15 - scaled x1000
$ws.Range("B15").Value = "51,305"

Set-TextValue "B16" "12.65"
Set-TextValue "B17" "157.85"
Set-TextValue "B18" "59.12"
Set-TextValue "B19" "12.65"
Set-TextValue "B20" "120.92"
Set-TextValue "B21" "103.73"
Set-TextValue "B22" "9.35"
Set-TextValue "B23" "50.3"

# Row 24 - scaled x1000
$ws.Range("B24").Value = "41,835"

Set-TextValue "B25" "125.11"
Set-TextValue "B26" "25.79"
Set-TextValue "B27" "19.06"

# Row 28 - scaled x1000
$ws.Range("B28").Value = "51,305"

Set-TextValue "B30" "74.8"

$helper.Clear()
$ws.Range("A1").Select()
